# expenses-ytd.xlsx maintenance edit:
#  - bump the report's "as of" date in the period label on the summary sheet
#  - shift a handful of expense dates on the expenses sheet by one day
#  - swap the custom "₪#,##0.00" currency number format used by the
#    amount column for the plain built-in "#,##0.00" format

$wb = $excel.ActiveWorkbook
$wsSummary  = $wb.Worksheets.Item("סיכום")
$wsExpenses = $wb.Worksheets.Item("הוצאות")

# --- Summary sheet: update the period label -------------------------------
$wsSummary.Range("A2").Value = "תקופה: 2026-01-01 עד 2026-02-13"

# --- Expenses sheet: nudge a few invoice dates forward one day ------------
# These cells hold plain text dates (not real Excel dates), so we enter
# them with a leading apostrophe to force text and then reapply the
# worksheet's default "Normal" style so no visible formatting changes.
function Set-TextDate($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextDate $wsExpenses.Range("A8")  "2026-01-18"
Set-TextDate $wsExpenses.Range("A9")  "2026-01-21"
Set-TextDate $wsExpenses.Range("A13") "2026-02-01"
Set-TextDate $wsExpenses.Range("A14") "2026-02-04"
Set-TextDate $wsExpenses.Range("A15") "2026-02-07"
Set-TextDate $wsExpenses.Range("A16") "2026-02-10"

# --- Expenses sheet: amount column number format ---------------------------
# Replace the custom currency format (₪#,##0.00) with the built-in
# #,##0.00 format.
$wsExpenses.Range("C2:C16").NumberFormat = "#,##0.00"
